$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.748.42"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "'1.731.66"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").Value = "'0.9976"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'242.51"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'0.9980"
$ws.Range("D7").Value = "'0.4932"
$ws.Range("E7").Value = "  +1.45%  "
$ws.Range("D8").Value = "'0.2625"
$ws.Range("E8").Value = "  +0.49%  "
$ws.Range("D9").Value = "'0.06221"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "'1.725.82"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "'15.97"
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("D12").Value = "'0.06996"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "'0.6127"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("D14").Value = "'4.507"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'77.31"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "'0.9978"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "'26.533.37"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "'0.9977"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "'0.000007230"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").Value = "'11.46"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'1.948.61"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'4.495"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'8.581"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'5.110"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Value = "'138.38"
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +0.60%  "
$ws.Range("D27").Value = "'1.778"
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").Value = "'106.57"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "'3.940"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "'3.680"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").Value = "'0.04484"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "'1.003"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").Value = "'0.6248"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "'0.9402"
$ws.Range("E37").Value = "  +3.59%  "
$ws.Range("D38").Value = "'2.051"
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("D39").Value = "'2.424"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "'0.9978"
$ws.Range("D41").Value = "'0.01515"
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("D42").Value = "'5.590"
$ws.Range("E42").Value = "  +3.54%  "
$ws.Range("D43").Value = "'99.44"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "'0.3866"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "'6.962"
$ws.Range("E45").Value = "  +4.14%  "
$ws.Range("D46").Value = "'0.1161"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "'0.05384"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").Value = "'7.879"
$ws.Range("E48").Value = "  +2.10%  "
$ws.Range("D49").Value = "'30.24"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").Value = "'51.85"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").Value = "'1.238"
$ws.Range("E51").Value = "  -1.07%  "
